$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("growCapacity")

# Row 2: Q128 declaration answer -> I_CANT_TELL (was PROBABLY_NOT)
$ws.Range("C2").Value = "I_CANT_TELL{}; "
$ws.Range("D2").Value = "I_CANT_TELL"

# Row 3: Q128-137 body answer -> PROBABLY_NOT (was I_CANT_TELL)
$ws.Range("C3").Value = "PROBABLY_NOT{}; "
$ws.Range("D3").Value = "PROBABLY_NOT"

# Row 4: Q129 invocation answer -> NO (was YES{found it})
$ws.Range("C4").Value = "NO{}; "
$ws.Range("D4").Value = "NO"

# Row 6: Q131 for-loop construct answer -> I_CANT_TELL (was PROBABLY_YES{...})
$ws.Range("C6").Value = "I_CANT_TELL{}; "
$ws.Range("D6").Value = "I_CANT_TELL"
$ws.Range("C6").WrapText = $false
